# Reorder the worksheet tabs so "review_info" comes before "hotel_info",
# and add a new "State" column to the hotel_info sheet (between Hotel_Name and City)
# populated with "Louisiana" for the existing hotel row.

$wb = $excel.ActiveWorkbook

# 1. Move "review_info" sheet so it appears before "hotel_info"
$wsReview = $wb.Worksheets.Item("review_info")
$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsReview.Move($wsHotel)

# NOTE: worksheet references captured before a Move() become stale (they track
# sheet *position*, not the logical sheet), so re-fetch "hotel_info" by name now.
$wsHotel = $wb.Worksheets.Item("hotel_info")

# 2. Insert a new column before the "City" column (column C) on hotel_info
#    so it lands right after "Hotel_Name" (column B).
$wsHotel.Columns.Item(3).Insert()

# 3. Populate the new column's header and value
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"
